$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")
$ws.Range("A1").Value = "Test"
